# ClientList.xlsx had picked up a pile of scratch/test sheets (Colton,
# Jim, Sheet1, Jeff) alongside the real "Client info" roster. This pass
# finishes the CreateClient / SelectVendor / SelectClient flow: a brand
# new client called "alexander" gets created, and the workbook is
# trimmed back down to just that client's sheet plus the client roster.

$wb = $excel.ActiveWorkbook

# The old "Client info" sheet carried a mailto hyperlink (with its own
# "Hyperlink" cell style) on Colton's row. Strip the hyperlink and the
# now-unused style before the sheet itself is discarded.
[void]$wb.Worksheets.Item("Client info").Hyperlinks.Delete()
[void]$wb.Styles.Item("Hyperlink").Delete()

# Drop the other scratch client sheets - only the new client and the
# roster sheet survive.
[void]$wb.Worksheets.Item("Colton").Delete()
[void]$wb.Worksheets.Item("Jeff").Delete()
[void]$wb.Worksheets.Item("Client info").Delete()

# "Jim" is repurposed as the newly created client sheet, "alexander".
$wsClient = $wb.Worksheets.Item("Jim")
$wsClient.Name = "alexander"
$wsClient.Range("A1").Value = "alexander"

# "Sheet1" becomes the trimmed "Client info" roster, holding just the
# new client's name/email pair.
$wsInfo = $wb.Worksheets.Item("Sheet1")
$wsInfo.Name = "Client info"
$wsInfo.Range("A1").Value = "alexander"
$wsInfo.Range("B1").Value = "alexander@alexander.com"

# The new client sheet is the one left active/selected.
$wsClient.Activate()
